$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.574.04"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.41%  "

# Row 3
$ws.Range("D3").Value = "'1.828.73"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'316.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.14%  "

# Row 6
$ws.Range("E6").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").Value = "'0.5333"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "'0.3982"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.32%  "

# Row 9
$ws.Range("D9").Value = "'0.07841"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.86%  "

# Row 10
$ws.Range("D10").Value = "'42.08"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.24%  "

# Row 11
$ws.Range("D11").Value = "'1.118"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.92%  "

# Row 12
$ws.Range("D12").Value = "'6.345"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.09%  "

# Row 13
$ws.Range("D13").Value = "'21.08"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.67%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.578"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.74%  "

# Row 15
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "'1.001"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.07%  "

# Row 16
$ws.Range("D16").Value = "'1.830.95"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.22%  "

# Row 17
$ws.Range("D17").Value = "'93.28"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.94%  "

# Row 18
$ws.Range("D18").Value = "'0.00001094"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.83%  "

# Row 19
$ws.Range("D19").Value = "'0.06562"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.79%  "

# Row 20
$ws.Range("D20").Value = "'17.82"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.30%  "

# Row 21
$ws.Range("E21").Value = "  +0.13%  "

# Row 22
$ws.Range("D22").Value = "'6.108"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.09%  "

# Row 23
$ws.Range("D23").Value = "'28.592.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.34%  "

# Row 24
$ws.Range("D24").Value = "'11.21"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.06%  "

# Row 25
$ws.Range("D25").Value = "'2.238"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.24%  "

# Row 26
$ws.Range("E26").Value = "  +1.74%  "

# Row 27
$ws.Range("D27").Value = "'157.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.52%  "

# Row 28
$ws.Range("D28").Value = "'2.039.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.18%  "

# Row 29
$ws.Range("D29").Value = "'2.416"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.06%  "

# Row 30
$ws.Range("D30").Value = "'125.47"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.75%  "

# Row 31
$ws.Range("E31").Value = "  +3.39%  "

# Row 32
$ws.Range("D32").Value = "'0.1123"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.71%  "

# Row 33
$ws.Range("D33").Value = "'5.746"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.98%  "

# Row 34
$ws.Range("E34").Value = "  +0.75%  "

# Row 35
$ws.Range("D35").Value = "'0.07320"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.05%  "

# Row 36
$ws.Range("D36").Value = "'0.2267"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.88%  "

# Row 37
$ws.Range("D37").Value = "'8.996"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.45%  "

# Row 38
$ws.Range("D38").Value = "'0.02348"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.08%  "

# Row 39
$ws.Range("D39").Value = "'5.213"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.66%  "

# Row 40
$ws.Range("E40").Value = "  +2.49%  "

# Row 41
$ws.Range("D41").Value = "'0.6297"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.99%  "

# Row 42
$ws.Range("D42").Value = "'1.198"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.38%  "

# Row 43
$ws.Range("E43").Value = "  +0.04%  "

# Row 44
$ws.Range("E44").Value = "  -3.26%  "

# Row 45
$ws.Range("D45").Value = "'13.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.41%  "

# Row 46
$ws.Range("D46").Value = "'0.5932"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.86%  "

# Row 47
$ws.Range("D47").Value = "'3.716"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.80%  "

# Row 48
$ws.Range("D48").Value = "'125.59"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.16%  "

# Row 49
$ws.Range("D49").Value = "'1.997"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.62%  "

# Row 50
$ws.Range("D50").Value = "'1.194"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.84%  "

# Row 51
$ws.Range("D51").Value = "'0.06952"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.92%  "
